# Update header labels on the "model" worksheet:
#   A1: "account"  -> "account_number"
#   C1: "weight"   -> "model_weight"
# and move the active selection from E24 to C1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

$ws.Range("A1").Value = "account_number"
$ws.Range("C1").Value = "model_weight"

$ws.Range("C1").Select()
